$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ------------------------------------------------------------------
# 1) Update the summary figures near the top of the account statement
# ------------------------------------------------------------------
$ws.Range("E11").Value = 332521   # VALOR MORA (total)
$ws.Range("C13").Value = 8        # Cant. Trabajadores
$ws.Range("F13").Value = 11       # Cant. Periodos

# ------------------------------------------------------------------
# 2) The worker/period table used to span rows 16:31 (16 rows).
#    The new data set only has 12 rows, so remove 4 rows from the
#    middle of the block (this shifts the bottom-bordered "last row"
#    of the table - plus the signature block below it - up by 4,
#    which is exactly the B2:J37 -> B2:J33 shrink seen in the diff)
#    while keeping that last row's special border styling intact.
# ------------------------------------------------------------------
$ws.Rows("16:19").Delete()

# ------------------------------------------------------------------
# 3) Write the refreshed worker/period rows (16:27)
# ------------------------------------------------------------------
$data = @(
    @("CC", "1050959757", "JEAN CARLOS CARRASCAL DIAZ", "1811", 55120, 877803),
    @("CC", "70440095", "JUAN CARLOS PEREZ RIOS", "2101", 7268, 908526),
    @("CC", "1128051813", "MARGOTH ELENA PORRAS SALCEDO", "2111", 36341, 908526),
    @("CC", "1043635331", "JEFERSON JOSE CANTILLO MEJIA", "2203", 40000, 1300000),
    @("CC", "1043635331", "JEFERSON JOSE CANTILLO MEJIA", "2204", 40000, 1300000),
    @("CC", "1043635331", "JEFERSON JOSE CANTILLO MEJIA", "2205", 40000, 1300000),
    @("CC", "1043635331", "JEFERSON JOSE CANTILLO MEJIA", "2206", 40000, 1300000),
    @("CC", "1043635331", "JEFERSON JOSE CANTILLO MEJIA", "2207", 40000, 1300000),
    @("CC", "1051827123", "JEAN CARLOS CASTELLAR HERRERA", "2302", 22027, 1652000),
    @("CC", "1047493632", "FABIO ANTONIO RODRIGUEZ AGUIRRE", "2302", 4058, 1521708),
    @("CC", "1047420154", "MARTHA PATRICIA BRAVO VILLERO", "2307", 4640, 1160000),
    @("CC", "1047473702", "CARLOS ANTONIO FERNANDEZ SALAMANCA", "2407", 3067, 2300000)
)

$row = 16
foreach ($r in $data) {
    $ws.Cells.Item($row, 2).Value = $r[0]
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    $ws.Cells.Item($row, 6).Value = $r[4]
    $ws.Cells.Item($row, 7).Value = $r[5]
    $row = $row + 1
}
